$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.114.80'
$ws.Range("E2").Value = '  -3.48%  '
$ws.Range("D3").Value = '3.340.72'
$ws.Range("E3").Value = '  -5.38%  '
$ws.Range("E4").Value = '  +0.10%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '551.08'
$r.Style = "Normal"
$ws.Range("E5").Value = '  -5.26%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '173.73'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -3.42%  '
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.613'
$r.Style = "Normal"
$ws.Range("E7").Value = '  -2.86%  '
$ws.Range("D8").Value = '3.334.83'
$ws.Range("E8").Value = '  -5.36%  '
$ws.Range("E9").Value = '  -0.01%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.625'
$r.Style = "Normal"
$ws.Range("E10").Value = '  -2.40%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.160'
$r.Style = "Normal"
$ws.Range("E11").Value = '  -1.65%  '
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '53.72'
$r.Style = "Normal"
$ws.Range("E12").Value = '  -3.80%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.0000272'
$r.Style = "Normal"
$ws.Range("E13").Value = '  -4.10%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '9.02'
$r.Style = "Normal"
$ws.Range("E14").Value = '  -3.21%  '
$ws.Range("D15").Value = '3.875.80'
$ws.Range("E15").Value = '  -5.31%  '
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '18.29'
$r.Style = "Normal"
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("E17").Value = '  -3.37%  '
$ws.Range("D18").Value = '3.338.01'
$ws.Range("E18").Value = '  -5.61%  '
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '11.74'
$r.Style = "Normal"
$ws.Range("E19").Value = '  -2.84%  '
$ws.Range("D20").Value = '64.113.88'
$ws.Range("E20").Value = '  -3.44%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '0.974'
$r.Style = "Normal"
$ws.Range("E21").Value = '  -3.62%  '
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '424.92'
$r.Style = "Normal"
$ws.Range("E22").Value = '  +2.08%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '4.82'
$r.Style = "Normal"
$ws.Range("E23").Value = '  +10.71%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '4.08'
$r.Style = "Normal"
$ws.Range("E24").Value = '  -4.46%  '
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '83.91'
$r.Style = "Normal"
$ws.Range("E25").Value = '  -2.23%  '
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '13.07'
$r.Style = "Normal"
$ws.Range("E26").Value = '  -1.91%  '
$ws.Range("E27").Value = '  -5.07%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '2.81'
$r.Style = "Normal"
$ws.Range("E28").Value = '  -2.45%  '
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '8.59'
$r.Style = "Normal"
$ws.Range("E29").Value = '  -6.22%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '29.61'
$r.Style = "Normal"
$ws.Range("E30").Value = '  -3.18%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '6.61'
$r.Style = "Normal"
$ws.Range("E31").Value = '  -0.15%  '
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '592.65'
$r.Style = "Normal"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '11.40'
$r.Style = "Normal"
$ws.Range("E33").Value = '  -3.18%  '
$ws.Range("E34").Value = '  -3.97%  '
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '58.11'
$r.Style = "Normal"
$ws.Range("E35").Value = '  -3.01%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("E37").Value = '  -9.08%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '3.50'
$r.Style = "Normal"
$ws.Range("E38").Value = '  -3.51%  '
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '35.48'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -5.20%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0752'
$ws.Range("E40").Value = '  -7.62%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.365'
$r.Style = "Normal"
$ws.Range("E41").Value = '  -5.74%  '
$ws.Range("D42").Value = '3.091.99'
$ws.Range("E42").Value = '  -4.98%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E43").Value = '  +0.15%  '
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '2.80'
$r.Style = "Normal"
$ws.Range("E44").Value = '  -5.01%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '3.21'
$r.Style = "Normal"
$ws.Range("E45").Value = '  -4.27%  '
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '0.0406'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -4.29%  '
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '2.45'
$r.Style = "Normal"
$ws.Range("E47").Value = '  -5.01%  '
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("E49").Value = '  -4.79%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '136.17'
$r.Style = "Normal"
$ws.Range("E50").Value = '  -2.64%  '
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '8.17'
$r.Style = "Normal"
$ws.Range("E51").Value = '  -5.88%  '
